$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing cells (NaN <-> numeric value swaps, plus a few pure numeric updates) ---
$ws.Range("W41").Value = "NaN"
$ws.Range("W42").Value = "NaN"
$ws.Range("W43").Value = "NaN"
$ws.Range("C44").Value = 1
$ws.Range("J51").Value = 2
$ws.Range("J52").Value = 2
$ws.Range("J53").Value = 2
$ws.Range("J54").Value = 2
$ws.Range("AA59").Value = 1
$ws.Range("J114").Value = 26
$ws.Range("J115").Value = 31
$ws.Range("AI123").Value = 8
$ws.Range("J156").Value = 2044
$ws.Range("J157").Value = 2144
$ws.Range("J158").Value = 2496
$ws.Range("J159").Value = 2620
$ws.Range("J160").Value = 2795
$ws.Range("J161").Value = 2945
$ws.Range("J162").Value = 3069
$ws.Range("J163").Value = 3330
$ws.Range("J164").Value = 3513
$ws.Range("J165").Value = 3752
$ws.Range("J166").Value = 3926
$ws.Range("J167").Value = 4141
$ws.Range("J168").Value = 4350
$ws.Range("J169").Value = 4440
$ws.Range("J170").Value = 4615
$ws.Range("J171").Value = 4724
$ws.Range("J172").Value = 4844
$ws.Range("J173").Value = 5077
$ws.Range("J174").Value = 5183
$ws.Range("J175").Value = 5210
$ws.Range("J176").Value = 5351
$ws.Range("J177").Value = 5540
$ws.Range("J178").Value = 5729
$ws.Range("J179").Value = 5866
$ws.Range("J180").Value = 6073
$ws.Range("J181").Value = 6171
$ws.Range("J182").Value = 6313
$ws.Range("J183").Value = 6322
$ws.Range("J184").Value = 6527
$ws.Range("J185").Value = 6631
$ws.Range("J186").Value = 6795
$ws.Range("J187").Value = 6903
$ws.Range("J188").Value = 7035
$ws.Range("J189").Value = 7114
$ws.Range("J190").Value = 7180
$ws.Range("J191").Value = 7240
$ws.Range("J192").Value = 7250
$ws.Range("J193").Value = 7309
$ws.Range("J194").Value = 7320
$ws.Range("J195").Value = 7347
$ws.Range("J196").Value = 7358
$ws.Range("J197").Value = 7388
$ws.Range("J198").Value = 7537
# --- Append new row 207 (next date in the Colombia COVID case series) ---
$row207Values = @(44101,813056,2742,112052,67320,264271,28589,6808,5664,8429,9223,19462,3990,23688,32779,7985,11036,15101,14761,18177,15507,3733,3328,10720,30364,13998,11976,61081,2210,1138,754,474,767,479,741,2062,5763,38032,9934,2565,47280,1106,22815,1527,10484,1674,1606,8340,2015,964,2501,2688,64947,14101,6595,9862,7204,257,1473,2731,744,2166,9848,9559,10681,14323,1967,904,13941,11360,13282,3110,2248,5981,4903,2378,6035,3799,2231,1072,3042,2244,2027,1849,6571,2216,1482,1818,2133,2238,2635,1757,1221,1224,1050,3454,1492,968,1128,1738,1624,814,900,1354,1689,1568,1640,1260,334,371,836,775,495,544,386,671,760,527,505,374,524,138457,343352,19985,149728,92471,45920,12929)
for ($col = 1; $col -le $row207Values.Length; $col++) {
    $ws.Cells.Item(207, $col).Value = $row207Values[$col - 1]
}

# --- Update the visible selection to the last cell touched, matching the saved view state ---
$ws.Range("DX207").Select()
